$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1250032.9
$ws.Range("I6").Value = 1250032.9
$ws.Range("K6").Value = 3750098.7
$ws.Range("M6").Value = -3749986.7

$ws.Range("H12").Value = 245.125
$ws.Range("I12").Value = 137.28572
$ws.Range("K12").Value = 137.28572
$ws.Range("M12").Value = 32.71428

$ws.Range("H13").Value = 14881.5
$ws.Range("J13").Value = 7857.8
$ws.Range("L13").Value = 7857.8
$ws.Range("N13").Value = -8195.799999999999

$ws.Range("H38").Value = 247.5
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 247.5
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 742.5
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -1486.5

$ws.Range("H43").Value = 9719.157999999999
$ws.Range("J43").Value = 5477
$ws.Range("L43").Value = 5477
$ws.Range("N43").Value = -5615

$ws.Range("H92").Value = 657.9048
$ws.Range("I92").Value = 657.9048
$ws.Range("K92").Value = 657.9048
$ws.Range("M92").Value = 590.0952

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws.Range("H107").Value = 968.25
$ws.Range("I107").Value = 1035.1428
$ws.Range("K107").Value = 1035.1428
$ws.Range("M107").Value = 884.8571999999999

$ws.Range("H135").Value = 1767.6666
$ws.Range("I135").Value = 1344
$ws.Range("J135").Value = 3250.5
$ws.Range("K135").Value = 12096
$ws.Range("L135").Value = 29254.5
$ws.Range("M135").Value = -9561
$ws.Range("N135").Value = -34324.5

$ws.Range("H138").Value = 6581995.5
$ws.Range("I138").Value = 1703.5
$ws.Range("J138").Value = 7579009.5
$ws.Range("K138").Value = 5110.5
$ws.Range("L138").Value = 22737028.5
$ws.Range("M138").Value = 29.5
$ws.Range("N138").Value = -22747308.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8366.834999999999
$ws.Range("I32").Value = 3562.127
$ws.Range("K32").Value = 3562.127
$ws.Range("M32").Value = -3275.127

$ws.Range("H132").Value = 2281.1355
$ws.Range("I132").Value = 2060.26
$ws.Range("J132").Value = 3508.2222
$ws.Range("K132").Value = 6180.780000000001
$ws.Range("L132").Value = 10524.6666
$ws.Range("M132").Value = -3650.780000000001
$ws.Range("N132").Value = -15584.6666

$ws.Range("H135").Value = 58500
$ws.Range("J135").Value = 58500
$ws.Range("L135").Value = 58500
$ws.Range("N135").Value = -68640

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 29757.5
$ws.Range("J81").Value = 29757.5
$ws.Range("L81").Value = 29757.5
$ws.Range("N81").Value = -31879.5

$ws.Range("H84").Value = 29757.5
$ws.Range("J84").Value = 29757.5
$ws.Range("L84").Value = 89272.5
$ws.Range("N84").Value = -99880.5

$ws.Range("H107").Value = 1202.9412
$ws.Range("I107").Value = 1088.3572
$ws.Range("J107").Value = 1737.6666
$ws.Range("K107").Value = 1088.3572
$ws.Range("L107").Value = 1737.6666
$ws.Range("M107").Value = 831.6428000000001
$ws.Range("N107").Value = -5577.6666

$ws.Range("H135").Value = 73097.5
$ws.Range("J135").Value = 73097.5
$ws.Range("L135").Value = 73097.5
$ws.Range("N135").Value = -83237.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6116.875
$ws.Range("I16").Value = 6759.4
$ws.Range("K16").Value = 6759.4
$ws.Range("M16").Value = -6472.4

$ws.Range("H23").Value = 48747.75
$ws.Range("J23").Value = 48747.75
$ws.Range("L23").Value = 48747.75
$ws.Range("N23").Value = -49227.75

$ws.Range("H27").Value = 48747.75
$ws.Range("J27").Value = 48747.75
$ws.Range("L27").Value = 48747.75
$ws.Range("N27").Value = -49131.75

$ws.Range("H31").Value = 69587.664
$ws.Range("I31").Value = 85670.5
$ws.Range("J31").Value = 5256.3335
$ws.Range("K31").Value = 85670.5
$ws.Range("L31").Value = 5256.3335
$ws.Range("M31").Value = -85375.5
$ws.Range("N31").Value = -5846.3335

$ws.Range("H34").Value = 69587.664
$ws.Range("I34").Value = 85670.5
$ws.Range("J34").Value = 5256.3335
$ws.Range("K34").Value = 85670.5
$ws.Range("L34").Value = 5256.3335
$ws.Range("M34").Value = -85468.5
$ws.Range("N34").Value = -5660.3335

$ws.Range("H41").Value = 19882.572
$ws.Range("J41").Value = 19882.572
$ws.Range("L41").Value = 19882.572
$ws.Range("N41").Value = -20738.572

$ws.Range("H58").Value = 3987.925
$ws.Range("I58").Value = 3770
$ws.Range("J58").Value = 4392.643
$ws.Range("K58").Value = 3770
$ws.Range("L58").Value = 4392.643
$ws.Range("M58").Value = -3567
$ws.Range("N58").Value = -4798.643

$ws.Range("H62").Value = 11333.333
$ws.Range("J62").Value = 14500
$ws.Range("L62").Value = 14500
$ws.Range("N62").Value = -15748

$ws.Range("H65").Value = 11333.333
$ws.Range("J65").Value = 14500
$ws.Range("L65").Value = 72500
$ws.Range("N65").Value = -78740

$ws.Range("H93").Value = 16726.25
$ws.Range("I93").Value = 9301.666999999999
$ws.Range("K93").Value = 9301.666999999999
$ws.Range("M93").Value = -7429.666999999999

$ws.Range("H103").Value = 36666.668
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

$ws.Range("H113").Value = 6116.875
$ws.Range("I113").Value = 6759.4
$ws.Range("K113").Value = 6759.4
$ws.Range("M113").Value = -4589.4

$ws.Range("H122").Value = 1437.8077
$ws.Range("I122").Value = 1021.8947
$ws.Range("J122").Value = 2566.7144
$ws.Range("K122").Value = 3065.6841
$ws.Range("L122").Value = 7700.1432
$ws.Range("M122").Value = -615.6840999999999
$ws.Range("N122").Value = -12600.1432

$ws.Range("H134").Value = 15409.904
$ws.Range("I134").Value = 9365.210999999999
$ws.Range("K134").Value = 28095.633
$ws.Range("M134").Value = -25560.633

$ws.Range("H136").Value = 3987.925
$ws.Range("I136").Value = 3770
$ws.Range("J136").Value = 4392.643
$ws.Range("K136").Value = 11310
$ws.Range("L136").Value = 13177.929
$ws.Range("M136").Value = -8760
$ws.Range("N136").Value = -18277.929

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 895.625
$ws.Range("I18").Value = 855.3333
$ws.Range("K18").Value = 2565.9999
$ws.Range("M18").Value = -2396.9999

$ws.Range("H92").Value = 1549.7
$ws.Range("J92").Value = 1728.875
$ws.Range("L92").Value = 5186.625
$ws.Range("N92").Value = -7682.625

$ws.Range("H122").Value = 1258.1666
$ws.Range("J122").Value = 1389.8
$ws.Range("L122").Value = 12508.2
$ws.Range("N122").Value = -17408.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 19000
$ws.Range("J27").Value = 19000
$ws.Range("L27").Value = 19000
$ws.Range("N27").Value = -19332

$ws.Range("H97").Value = 1234.0526
$ws.Range("I97").Value = 1008.4167
$ws.Range("K97").Value = 1008.4167
$ws.Range("M97").Value = -512.4167

$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164

$ws.Range("H102").Value = 41667456
$ws.Range("I102").Value = 821.5
$ws.Range("K102").Value = 821.5
$ws.Range("M102").Value = 800.5

$ws.Range("H132").Value = 2606.976
$ws.Range("I132").Value = 2309.125
$ws.Range("K132").Value = 6927.375
$ws.Range("M132").Value = -4397.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H40").Value = 2558.5
$ws.Range("I40").Value = 1711.9286
$ws.Range("K40").Value = 1711.9286
$ws.Range("M40").Value = -1575.9286

$ws.Range("H68").Value = 1669166.6
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 10000000
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 10000000
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -10001498

$ws.Range("H71").Value = 1669166.6
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 10000000
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 50000000
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -50007488

$ws.Range("H132").Value = 4115.6665
$ws.Range("I132").Value = 4237.4614
$ws.Range("J132").Value = 3799
$ws.Range("K132").Value = 12712.3842
$ws.Range("L132").Value = 11397
$ws.Range("M132").Value = -10182.3842
$ws.Range("N132").Value = -16457

$ws.Range("H136").Value = 3205.2
$ws.Range("I136").Value = 2791.2856
$ws.Range("K136").Value = 8373.856800000001
$ws.Range("M136").Value = -5823.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7499.2964
$ws.Range("J81").Value = 4578.6313
$ws.Range("L81").Value = 9157.2626
$ws.Range("N81").Value = -11279.2626

$ws.Range("H84").Value = 7499.2964
$ws.Range("J84").Value = 4578.6313
$ws.Range("L84").Value = 45786.313
$ws.Range("N84").Value = -56394.313

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H126").Value = 2322.4614
$ws.Range("I126").Value = 2155
$ws.Range("K126").Value = 6465
$ws.Range("M126").Value = -3995

$ws.Range("H129").Value = 72500
$ws.Range("J129").Value = 72500
$ws.Range("L129").Value = 72500
$ws.Range("N129").Value = -72500
